$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("A4").Value = 131256691
$ws.Range("B4").Value = 57884
$ws.Range("E4").Value = 100109
$ws.Range("F4").Value = "Tretåig hackspett"
$ws.Range("G4").Value = "Picoides tridactylus"
$ws.Range("H4").Value = "(Linnaeus, 1758)"
$ws.Range("M4").Value = "äldre spår"
$ws.Range("Q4").Value = 488667
$ws.Range("R4").Value = 6665262
$ws.Range("Z4").Value = "10:55"
$ws.Range("AB4").Value = "10:55"
$ws.Range("AC4").Value = "Ringhack på gran."

# Row 5
$ws.Range("A5").Value = 131260583
$ws.Range("B5").Value = 57884
$ws.Range("E5").Value = 100109
$ws.Range("F5").Value = "Tretåig hackspett"
$ws.Range("G5").Value = "Picoides tridactylus"
$ws.Range("H5").Value = "(Linnaeus, 1758)"
$ws.Range("M5").Value = "färska spår"
$ws.Range("Q5").Value = 488834
$ws.Range("R5").Value = 6665228
$ws.Range("Z5").Value = "15:30"
$ws.Range("AB5").Value = "15:30"
$ws.Range("AC5").Value = "Ringhack på tall."

# Row 6
$ws.Range("A6").Value = 131257424
$ws.Range("B6").Value = 79245
$ws.Range("E6").Value = 6425
$ws.Range("F6").Value = "Garnlav"
$ws.Range("G6").Value = "Alectoria sarmentosa"
$ws.Range("H6").Value = "(Ach.) Ach."
$ws.Range("M6").Value = ""
$ws.Range("Q6").Value = 488876
$ws.Range("R6").Value = 6665177
$ws.Range("Z6").Value = "11:33"
$ws.Range("AB6").Value = "11:33"
$ws.Range("AC6").Value = "Gran"

# Row 7
$ws.Range("A7").Value = 131255793
$ws.Range("B7").Value = 91833
$ws.Range("E7").Value = 5432
$ws.Range("F7").Value = "Granticka"
$ws.Range("G7").Value = "Porodaedalea chrysoloma s.lat."
$ws.Range("H7").Value = ""
$ws.Range("M7").Value = ""
$ws.Range("Q7").Value = 488817
$ws.Range("R7").Value = 6665110
$ws.Range("Z7").Value = "09:56"
$ws.Range("AB7").Value = "09:56"
$ws.Range("AC7").Value = "Flera fruktkroppar."

# Row 8
$ws.Range("A8").Value = 131256423
$ws.Range("B8").Value = 57881
$ws.Range("E8").Value = 100049
$ws.Range("F8").Value = "Spillkråka"
$ws.Range("G8").Value = "Dryocopus martius"
$ws.Range("H8").Value = "(Linnaeus, 1758)"
$ws.Range("I8").Value = "1"
$ws.Range("K8").Value = "adult"
$ws.Range("M8").Value = "spel/sång"
$ws.Range("Q8").Value = 488671
$ws.Range("R8").Value = 6665267
$ws.Range("S8").Value = 10
$ws.Range("Z8").Value = "10:40"
$ws.Range("AB8").Value = "10:40"
$ws.Range("AC8").Value = ""

# Row 9
$ws.Range("A9").Value = 131257188
$ws.Range("B9").Value = 91833
$ws.Range("E9").Value = 5432
$ws.Range("F9").Value = "Granticka"
$ws.Range("G9").Value = "Porodaedalea chrysoloma s.lat."
$ws.Range("H9").Value = ""
$ws.Range("I9").Value = ""
$ws.Range("K9").Value = ""
$ws.Range("M9").Value = ""
$ws.Range("Q9").Value = 488804
$ws.Range("R9").Value = 6665288
$ws.Range("S9").Value = 5
$ws.Range("Z9").Value = "11:17"
$ws.Range("AB9").Value = "11:17"
$ws.Range("AC9").Value = "Rikligt."

# Row 16
$ws.Range("B16").Value = 91833

# Row 21
$ws.Range("A21").Value = 131257316
$ws.Range("B21").Value = 79245
$ws.Range("E21").Value = 6425
$ws.Range("F21").Value = "Garnlav"
$ws.Range("G21").Value = "Alectoria sarmentosa"
$ws.Range("H21").Value = "(Ach.) Ach."
$ws.Range("Q21").Value = 488852
$ws.Range("R21").Value = 6665209
$ws.Range("Z21").Value = "11:28"
$ws.Range("AB21").Value = "11:28"
$ws.Range("AC21").Value = "Gran"

# Row 22
$ws.Range("A22").Value = 131255771
$ws.Range("B22").Value = 81230
$ws.Range("E22").Value = 1049
$ws.Range("F22").Value = "Kortskaftad ärgspik"
$ws.Range("G22").Value = "Microcalicium ahlneri"
$ws.Range("H22").Value = "Tibell"
$ws.Range("Q22").Value = 488818
$ws.Range("R22").Value = 6665110
$ws.Range("Z22").Value = "09:54"
$ws.Range("AB22").Value = "09:54"
$ws.Range("AC22").Value = ""

# Row 27
$ws.Range("A27").Value = 131255910
$ws.Range("Q27").Value = 488763
$ws.Range("R27").Value = 6665157
$ws.Range("Z27").Value = "10:03"
$ws.Range("AB27").Value = "10:03"
$ws.Range("AC27").Value = "Tall."

# Row 28
$ws.Range("A28").Value = 131258531
$ws.Range("B28").Value = 79245
$ws.Range("E28").Value = 6425
$ws.Range("F28").Value = "Garnlav"
$ws.Range("G28").Value = "Alectoria sarmentosa"
$ws.Range("H28").Value = "(Ach.) Ach."
$ws.Range("M28").Value = ""
$ws.Range("Q28").Value = 488725
$ws.Range("R28").Value = 6665212
$ws.Range("Z28").Value = "13:02"
$ws.Range("AB28").Value = "13:02"
$ws.Range("AC28").Value = "Gran"

# Row 29
$ws.Range("A29").Value = 131257239
$ws.Range("B29").Value = 57884
$ws.Range("E29").Value = 100109
$ws.Range("F29").Value = "Tretåig hackspett"
$ws.Range("G29").Value = "Picoides tridactylus"
$ws.Range("H29").Value = "(Linnaeus, 1758)"
$ws.Range("M29").Value = "färska spår"
$ws.Range("Q29").Value = 488852
$ws.Range("R29").Value = 6665286
$ws.Range("Z29").Value = "11:23"
$ws.Range("AB29").Value = "11:23"
$ws.Range("AC29").Value = "Barkfläk, hagelsalva."

# Row 35
$ws.Range("B35").Value = 91833
